$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "SC 92" row (28) first, then the "RM 232" row (26).
# Deleting the lower row first keeps the higher row's index stable.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# After the two deletions the rows 26-33 now hold (in order):
# SC 5, SC 101, SC 105, SC 119, SC 120, SC 132, SC 193, SC 232
# Apply the remaining value corrections for those rows.
$ws.Range("D27").Value = -14.6
$ws.Range("D28").ClearContents()
$ws.Range("D29").ClearContents()
$ws.Range("D30").Value = -13.6
$ws.Range("E30").Value = -5.7
$ws.Range("F30").Value = 16.89
$ws.Range("D32").ClearContents()

# Corrections to rows 2-25 (imputed / cleared values).
$ws.Range("F2").ClearContents()
$ws.Range("E6").Value = -5.7
$ws.Range("E8").ClearContents()
$ws.Range("E18").Value = -8.5
$ws.Range("E20").ClearContents()
$ws.Range("E23").Value = -7
$ws.Range("E25").ClearContents()
